$d = $word.ActiveDocument

$replacements = @(
    @{ old = "405÷4=101, 1"; new = "458÷3=152, 2" },
    @{ old = "353÷4=88, 1";  new = "141÷3=47, 0" },
    @{ old = "772÷7=110, 2"; new = "137÷2=68, 1" },
    @{ old = "457÷5=91, 2";  new = "329÷6=54, 5" },
    @{ old = "357÷4=89, 1";  new = "764÷2=382, 0" },
    @{ old = "268÷9=29, 7";  new = "750÷2=375, 0" },
    @{ old = "243÷6=40, 3";  new = "611÷5=122, 1" },
    @{ old = "851÷8=106, 3"; new = "845÷6=140, 5" },
    @{ old = "249÷6=41, 3";  new = "722÷5=144, 2" },
    @{ old = "605÷7=86, 3";  new = "222÷2=111, 0" },
    @{ old = "765÷2=382, 1"; new = "995÷6=165, 5" },
    @{ old = "635÷6=105, 5"; new = "999÷5=199, 4" },
    @{ old = "876÷3=292, 0"; new = "370÷4=92, 2" },
    @{ old = "871÷6=145, 1"; new = "884÷9=98, 2" },
    @{ old = "802÷7=114, 4"; new = "236÷5=47, 1" },
    @{ old = "565÷2=282, 1"; new = "782÷7=111, 5" },
    @{ old = "818÷5=163, 3"; new = "524÷8=65, 4" },
    @{ old = "773÷5=154, 3"; new = "809÷5=161, 4" },
    @{ old = "694÷7=99, 1";  new = "276÷4=69, 0" },
    @{ old = "191÷9=21, 2";  new = "514÷2=257, 0" },
    @{ old = "972÷3=324, 0"; new = "135÷5=27, 0" },
    @{ old = "940÷5=188, 0"; new = "314÷5=62, 4" },
    @{ old = "341÷6=56, 5";  new = "451÷4=112, 3" },
    @{ old = "979÷7=139, 6"; new = "608÷7=86, 6" },
    @{ old = "879÷8=109, 7"; new = "319÷9=35, 4" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
